$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3009
    $ws.Range("F4").Value = 105
    $ws.Range("F5").Value = 6764
    $ws.Range("F6").Value = 1755
    $ws.Range("F10").Value = 126
    $ws.Range("F11").Value = 133
}
